$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2:AJ2").ClearContents()

# Row 3
$ws.Range("D3:AF3").ClearContents()
$ws.Range("AI3").ClearContents()
$ws.Range("AG3").Value = 200
$ws.Range("AH3").Value = 0.65
$ws.Range("AJ3").Value = 9100836

# Row 4
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()
$ws.Range("V4").ClearContents()
$ws.Range("Y4:Z4").ClearContents()
$ws.Range("D4").Value = 1536
$ws.Range("E4").Value = 94
$ws.Range("F4").Value = 94
$ws.Range("G4").Value = 83
$ws.Range("H4").Value = 64
$ws.Range("I4").Value = 64
$ws.Range("K4").Value = 2103
$ws.Range("L4").Value = 229
$ws.Range("M4").Value = 1874
$ws.Range("N4").Value = 1874
$ws.Range("P4").Value = 455
$ws.Range("Q4").Value = 81
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = -18
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 81
$ws.Range("W4").Value = 6.14
$ws.Range("X4").Value = 4.18
$ws.Range("AA4").Value = 12.2
$ws.Range("AB4").Value = 332.28
$ws.Range("AC4").Value = 705
$ws.Range("AD4").Value = 23.32
$ws.Range("AE4").Value = 20604
$ws.Range("AF4").Value = 0.8
$ws.Range("AG4").Value = 400
$ws.Range("AH4").Value = 2.43
$ws.Range("AI4").Value = 56.68
$ws.Range("AJ4").Value = 9100836

# Row 5
$ws.Range("D5").Value = 1452
$ws.Range("E5").Value = 87
$ws.Range("F5").Value = 87
$ws.Range("G5").Value = 151
$ws.Range("H5").Value = 120
$ws.Range("I5").Value = 121
$ws.Range("J5").Value = -1
$ws.Range("K5").Value = 1985
$ws.Range("L5").Value = 202
$ws.Range("M5").Value = 1784
$ws.Range("N5").Value = 1781
$ws.Range("O5").Value = 3
$ws.Range("P5").Value = 455
$ws.Range("Q5").Value = 93
$ws.Range("R5").Value = -36
$ws.Range("S5").Value = -76
$ws.Range("T5").Value = 18
$ws.Range("U5").Value = 75
$ws.Range("V5").Value = 0
$ws.Range("W5").Value = 6.02
$ws.Range("X5").Value = 8.279999999999999
$ws.Range("Y5").Value = 6.62
$ws.Range("Z5").Value = 5.88
$ws.Range("AA5").Value = 11.31
$ws.Range("AB5").Value = 349.98
$ws.Range("AC5").Value = 1330
$ws.Range("AD5").Value = 10.98
$ws.Range("AE5").Value = 20246
$ws.Range("AF5").Value = 0.72
$ws.Range("AG5").Value = 500
$ws.Range("AH5").Value = 3.42
$ws.Range("AI5").Value = 36.35
$ws.Range("AJ5").Value = 9100836

# Row 6
$ws.Range("D6").Value = 1654
$ws.Range("E6").Value = 79
$ws.Range("F6").Value = 79
$ws.Range("G6").Value = 70
$ws.Range("H6").Value = 37
$ws.Range("I6").Value = 45
$ws.Range("K6").Value = 1888
$ws.Range("L6").Value = 205
$ws.Range("M6").Value = 1683
$ws.Range("N6").Value = 1676
$ws.Range("P6").Value = 455
$ws.Range("Q6").Value = -1
$ws.Range("R6").Value = -24
$ws.Range("S6").Value = -49
$ws.Range("T6").Value = 28
$ws.Range("U6").Value = -29
$ws.Range("V6").Value = 0
$ws.Range("W6").Value = 4.77
$ws.Range("X6").Value = 2.26
$ws.Range("Y6").Value = 2.58
$ws.Range("Z6").Value = 1.93
$ws.Range("AA6").Value = 12.2
$ws.Range("AB6").Value = 349.56
$ws.Range("AC6").Value = 490
$ws.Range("AD6").Value = 32.64
$ws.Range("AE6").Value = 19293
$ws.Range("AF6").Value = 0.83
$ws.Range("AG6").Value = 500
$ws.Range("AH6").Value = 3.13
$ws.Range("AI6").Value = 97.36
$ws.Range("AJ6").Value = 9100836

# Row 7
$ws.Range("U7").ClearContents()
$ws.Range("D7").Value = 1921
$ws.Range("E7").Value = 102
$ws.Range("G7").Value = 163
$ws.Range("H7").Value = 127
$ws.Range("I7").Value = 130
$ws.Range("K7").Value = 2360
$ws.Range("L7").Value = 570
$ws.Range("M7").Value = 1790
$ws.Range("N7").Value = 1790
$ws.Range("P7").Value = 460
$ws.Range("Q7").Value = -60
$ws.Range("R7").Value = 150
$ws.Range("S7").Value = -50
$ws.Range("T7").Value = 10
$ws.Range("W7").Value = 5.31
$ws.Range("X7").Value = 6.61
$ws.Range("Y7").Value = 7.5
$ws.Range("Z7").Value = 5.98
$ws.Range("AA7").Value = 31.84
$ws.Range("AC7").Value = 1428
$ws.Range("AD7").Value = 7.28
$ws.Range("AE7").Value = 20606
$ws.Range("AF7").Value = 0.5
$ws.Range("AG7").Value = 500
$ws.Range("AH7").Value = 4.81
$ws.Range("AI7").Value = 35

# Row 8
$ws.Range("U8").ClearContents()
$ws.Range("D8").Value = 1975
$ws.Range("E8").Value = 111
$ws.Range("G8").Value = 170
$ws.Range("H8").Value = 133
$ws.Range("I8").Value = 140
$ws.Range("K8").Value = 2450
$ws.Range("L8").Value = 570
$ws.Range("M8").Value = 1880
$ws.Range("N8").Value = 1880
$ws.Range("P8").Value = 460
$ws.Range("Q8").Value = 80
$ws.Range("R8").Value = 40
$ws.Range("S8").Value = -50
$ws.Range("T8").Value = 10
$ws.Range("W8").Value = 5.62
$ws.Range("X8").Value = 6.73
$ws.Range("Y8").Value = 7.63
$ws.Range("Z8").Value = 5.53
$ws.Range("AA8").Value = 30.32
$ws.Range("AC8").Value = 1538
$ws.Range("AD8").Value = 6.76
$ws.Range("AE8").Value = 21642
$ws.Range("AF8").Value = 0.48
$ws.Range("AG8").Value = 500
$ws.Range("AH8").Value = 4.81
$ws.Range("AI8").Value = 32.5

# Row 9
$ws.Range("U9").ClearContents()
$ws.Range("D9").Value = 2020
$ws.Range("E9").Value = 120
$ws.Range("G9").Value = 180
$ws.Range("H9").Value = 140
$ws.Range("I9").Value = 150
$ws.Range("K9").Value = 2540
$ws.Range("L9").Value = 570
$ws.Range("M9").Value = 1970
$ws.Range("N9").Value = 1980
$ws.Range("P9").Value = 460
$ws.Range("Q9").Value = 90
$ws.Range("R9").Value = 40
$ws.Range("S9").Value = -50
$ws.Range("T9").Value = 10
$ws.Range("W9").Value = 5.94
$ws.Range("X9").Value = 6.93
$ws.Range("Y9").Value = 7.77
$ws.Range("Z9").Value = 5.61
$ws.Range("AA9").Value = 28.93
$ws.Range("AC9").Value = 1648
$ws.Range("AD9").Value = 6.31
$ws.Range("AE9").Value = 22793
$ws.Range("AF9").Value = 0.46
$ws.Range("AG9").Value = 500
$ws.Range("AH9").Value = 4.81
$ws.Range("AI9").Value = 30.34
